# AYTO VIP S2 - add "Episode 18" results worksheet (after "Episode 17")
$wb = $excel.ActiveWorkbook

# --- Create the new sheet, positioned after the current last sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "Episode 18"

# --- Column headers (men), row 1, columns B..L ---
$menNames = @("Amadu","Calvin","Fabio","Luca","Lukas","Martin","Maurice","Max","Michael","Pharrell","Felix")
$menCols  = @("B","C","D","E","F","G","H","I","J","K","L")
for ($i = 0; $i -lt $menNames.Count; $i++) {
    $ws.Range($menCols[$i] + "1").Value = $menNames[$i]
}

# --- Row headers (women), column A, rows 2..11 ---
$womenNames = @("Anna","Cecilia","Celina","Franziska","Gina","Isabelle","Karina","Luisa","Ricarda","Zoe")
for ($i = 0; $i -lt $womenNames.Count; $i++) {
    $ws.Range("A" + ($i + 2)).Value = $womenNames[$i]
}

# Header formatting shared by the men-name row and the women-name column:
# bold text, thin box border on all sides, centered horizontally, top-aligned vertically
$headerRow = $ws.Range("B1:L1")
$headerRow.Font.Bold = $true
$headerRow.Borders.LineStyle = 1
$headerRow.HorizontalAlignment = -4108
$headerRow.VerticalAlignment = -4160

$headerCol = $ws.Range("A2:A11")
$headerCol.Font.Bold = $true
$headerCol.Borders.LineStyle = 1
$headerCol.HorizontalAlignment = -4108
$headerCol.VerticalAlignment = -4160

# --- Results grid B2:L11 ---
# Default every pairing to "no match" (0) with a light-gray fill
$grid = $ws.Range("B2:L11")
$grid.Value = 0
$grid.Interior.Color = 13882323

# Overwrite the confirmed "perfect match" pairings with 1, a magenta fill and
# light (near-white) font color so the value still reads against the fill
$matchCells = @("J2","B3","D4","C5","G6","I7","L7","E8","F9","H10","K11")
foreach ($ref in $matchCells) {
    $cell = $ws.Range($ref)
    $cell.Value = 1
    $cell.Interior.ColorIndex = 7
    $cell.Font.Color = 15856113
}

Write-Output "Added 'Episode 18' worksheet with match results"
